# Applies the "memory optimized code" pipeline update to SoIB_summaries.xlsx
#
# Summary of changes:
#  1. Sheet "Trends Status": Insufficient Data row (B8/C8) 308 -> 309
#  2. Sheet "Priority Status": High/Moderate/Low counts updated
#  3. Sheet "Species qualification": label + count updated
#  4. Sheet "High Priority break-up" renamed to "Interannual update - High Pri"
#     and its content replaced with new interannual breakdown (Trend New + IUCN)
#  5. New sheet "Major update - High Priority " added right after it, holding
#     the original "High Priority break-up" data (the old IUCN-only breakdown)
#
# NOTE: worksheet object references captured in variables can become stale
# once the sheet collection is restructured (Add/Move/rename shuffle the
# underlying index they resolve against), so sheets are re-fetched by name
# via $wb.Worksheets.Item(...) immediately before each use.

$wb = $excel.ActiveWorkbook

# --- 1. Trends Status: Insufficient Data row ---
$wb.Worksheets.Item("Trends Status").Range("B8").Value = 309
$wb.Worksheets.Item("Trends Status").Range("C8").Value = 309

# --- 2. Priority Status: updated counts ---
$wb.Worksheets.Item("Priority Status").Range("B2").Value = 103
$wb.Worksheets.Item("Priority Status").Range("B3").Value = 286
$wb.Worksheets.Item("Priority Status").Range("B4").Value = 554

# --- 3. Species qualification: label + count ---
$wb.Worksheets.Item("Species qualification").Range("A2").Value = "SoIB Assessment"
$wb.Worksheets.Item("Species qualification").Range("B2").Value = 309

# --- 5a. Create the new "Major update - High Priority " sheet, carrying over
#         the original "High Priority break-up" values, and position it right
#         after "High Priority break-up" ---
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "Major update - High Priority "

$wb.Worksheets.Item("Major update - High Priority ").Move($null, $wb.Worksheets.Item("High Priority break-up"))

$wsMajor = $wb.Worksheets.Item("Major update - High Priority ")
$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"
$wsMajor.Range("A1:E1").Font.Bold = $true
$wsMajor.Range("A1:E1").HorizontalAlignment = -4108

$wsMajor.Range("A2").Value = "IUCN"
$wsMajor.Range("B2").Value = 4
$wsMajor.Range("C2").Value = 100
$wsMajor.Range("D2").Value = 4
$wsMajor.Range("E2").Value = 100

# --- 5b. Rename "High Priority break-up" and replace its content with the
#         new interannual breakdown ---
$wb.Worksheets.Item("High Priority break-up").Name = "Interannual update - High Pri"

$wsInter = $wb.Worksheets.Item("Interannual update - High Pri")

$wsInter.Range("A1").Value = "Break-up"
$wsInter.Range("B1").Value = "High Species (no.)"
$wsInter.Range("C1").Value = "High Species (perc.)"
$wsInter.Range("D1").Value = "New High Species (no.)"
$wsInter.Range("E1").Value = "New High Species (perc.)"

$wsInter.Range("A2").Value = "Trend New"
$wsInter.Range("B2").Value = 97
$wsInter.Range("C2").Value = 94.2
$wsInter.Range("D2").Value = 97
$wsInter.Range("E2").Value = 98

$wsInter.Range("A3").Value = "IUCN"
$wsInter.Range("B3").Value = 6
$wsInter.Range("C3").Value = 5.8
$wsInter.Range("D3").Value = 2
$wsInter.Range("E3").Value = 2
